$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "1 Room"
$ws.Range("G3").Value = "1 Room"
$ws.Range("G8").Select()
